$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1585.8334
$ws.Range("I15").Value = 1585.8334
$ws.Range("K15").Value = 4757.5002
$ws.Range("M15").Value = -4588.5002
# Row 124
$ws.Range("H124").Value = 172666.67
$ws.Range("J124").Value = 172666.67
$ws.Range("L124").Value = 172666.67
$ws.Range("N124").Value = -182486.67
# Row 132
$ws.Range("H132").Value = 1616.1538
$ws.Range("I132").Value = 1454.375
$ws.Range("J132").Value = 1875
$ws.Range("K132").Value = 4363.125
$ws.Range("L132").Value = 5625
$ws.Range("M132").Value = -1833.125
$ws.Range("N132").Value = -10685
# Row 135
$ws.Range("H135").Value = 1713.875
$ws.Range("I135").Value = 1713.875
$ws.Range("K135").Value = 15424.875
$ws.Range("M135").Value = -12889.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 711.8421
$ws.Range("I2").Value = 757.7143
$ws.Range("J2").Value = 583.4
$ws.Range("K2").Value = 757.7143
$ws.Range("L2").Value = 583.4
$ws.Range("M2").Value = -644.7143
$ws.Range("N2").Value = -809.4
# Row 32
$ws.Range("H32").Value = 3221.25
$ws.Range("I32").Value = 2738.3076
$ws.Range("K32").Value = 2738.3076
$ws.Range("M32").Value = -2451.3076
# Row 61
$ws.Range("H61").Value = 1774.6666
$ws.Range("I61").Value = 1488.7333
$ws.Range("K61").Value = 1488.7333
$ws.Range("M61").Value = -1276.7333
# Row 74
$ws.Range("H74").Value = 24993622
$ws.Range("I74").Value = 66641330
$ws.Range("J74").Value = 4999.6
$ws.Range("K74").Value = 66641330
$ws.Range("L74").Value = 4999.6
$ws.Range("M74").Value = -66640456
$ws.Range("N74").Value = -6747.6
# Row 77
$ws.Range("H77").Value = 24993622
$ws.Range("I77").Value = 66641330
$ws.Range("J77").Value = 4999.6
$ws.Range("K77").Value = 333206650
$ws.Range("L77").Value = 24998
$ws.Range("M77").Value = -333202282
$ws.Range("N77").Value = -33734
# Row 116
$ws.Range("H116").Value = 711.8421
$ws.Range("I116").Value = 757.7143
$ws.Range("J116").Value = 583.4
$ws.Range("K116").Value = 757.7143
$ws.Range("L116").Value = 583.4
$ws.Range("M116").Value = 1536.2857
$ws.Range("N116").Value = -5171.4
# Row 132
$ws.Range("H132").Value = 2099.1538
$ws.Range("I132").Value = 994.875
$ws.Range("K132").Value = 2984.625
$ws.Range("M132").Value = -454.625
# Row 136
$ws.Range("H136").Value = 1774.6666
$ws.Range("I136").Value = 1488.7333
$ws.Range("K136").Value = 4466.199900000001
$ws.Range("M136").Value = -1916.199900000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 711.8421
$ws.Range("I3").Value = 757.7143
$ws.Range("J3").Value = 583.4
$ws.Range("K3").Value = 757.7143
$ws.Range("L3").Value = 583.4
$ws.Range("M3").Value = -643.7143
$ws.Range("N3").Value = -811.4
# Row 10
$ws.Range("H10").Value = 7006
$ws.Range("J10").Value = 7006
$ws.Range("L10").Value = 7006
$ws.Range("N10").Value = -7286
# Row 24
$ws.Range("H24").Value = 1208
$ws.Range("I24").Value = 1208
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 1208
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -973
$ws.Range("N24").ClearContents()
# Row 86
$ws.Range("H86").Value = 3511.182
$ws.Range("J86").Value = 4033
$ws.Range("L86").Value = 4033
$ws.Range("N86").Value = -6279
# Row 89
$ws.Range("H89").Value = 3511.182
$ws.Range("J89").Value = 4033
$ws.Range("L89").Value = 20165
$ws.Range("N89").Value = -31397

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4761.778
$ws.Range("J31").Value = 23998
$ws.Range("L31").Value = 23998
$ws.Range("N31").Value = -24588
# Row 34
$ws.Range("H34").Value = 4761.778
$ws.Range("J34").Value = 23998
$ws.Range("L34").Value = 23998
$ws.Range("N34").Value = -24402
# Row 134
$ws.Range("H134").Value = 1955.3334
$ws.Range("I134").Value = 2049.75
$ws.Range("K134").Value = 6149.25
$ws.Range("M134").Value = -3614.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 203512.8
$ws.Range("I80").Value = 4391.6665
$ws.Range("K80").Value = 13174.9995
$ws.Range("M80").Value = -12238.9995
# Row 83
$ws.Range("H83").Value = 203512.8
$ws.Range("I83").Value = 4391.6665
$ws.Range("K83").Value = 39524.9985
$ws.Range("M83").Value = -34844.9985
# Row 98
$ws.Range("H98").Value = 19212.428
$ws.Range("I98").Value = 5747.5
$ws.Range("J98").Value = 24598.4
$ws.Range("K98").Value = 17242.5
$ws.Range("L98").Value = 73795.20000000001
$ws.Range("M98").Value = -15744.5
$ws.Range("N98").Value = -76791.20000000001
# Row 131
$ws.Range("H131").Value = 1391.8
$ws.Range("J131").Value = 1589
$ws.Range("L131").Value = 4767
$ws.Range("N131").Value = -14847
# Row 132
$ws.Range("H132").Value = 3287.8333
$ws.Range("I132").Value = 2826.5
$ws.Range("J132").Value = 3749.1667
$ws.Range("K132").Value = 25438.5
$ws.Range("L132").Value = 33742.5003
$ws.Range("M132").Value = -22908.5
$ws.Range("N132").Value = -38802.5003

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 530.5
$ws.Range("J22").Value = 961
$ws.Range("L22").Value = 961
$ws.Range("N22").Value = -1551
# Row 27
$ws.Range("H27").Value = 530.5
$ws.Range("J27").Value = 961
$ws.Range("L27").Value = 961
$ws.Range("N27").Value = -1175
# Row 43
$ws.Range("H43").Value = 165713.14
$ws.Range("J43").Value = 165713.14
$ws.Range("L43").Value = 165713.14
$ws.Range("N43").Value = -166099.14
# Row 55
$ws.Range("H55").Value = 213.78572
$ws.Range("I55").Value = 199.5
$ws.Range("J55").Value = 299.5
$ws.Range("K55").Value = 199.5
$ws.Range("L55").Value = 299.5
$ws.Range("M55").Value = -26.5
$ws.Range("N55").Value = -645.5
# Row 61
$ws.Range("H61").Value = 4569.4287
$ws.Range("I61").Value = 4498.5
$ws.Range("K61").Value = 4498.5
$ws.Range("M61").Value = -4296.5
# Row 93
$ws.Range("H93").Value = 2730.1
$ws.Range("I93").Value = 2689.5715
$ws.Range("K93").Value = 2689.5715
$ws.Range("M93").Value = -1441.5715
# Row 113
$ws.Range("H113").Value = 4569.4287
$ws.Range("I113").Value = 4498.5
$ws.Range("K113").Value = 4498.5
$ws.Range("M113").Value = -2328.5
# Row 136
$ws.Range("H136").Value = 7999629.5
$ws.Range("I136").Value = 11426685
$ws.Range("J136").Value = 3166.3333
$ws.Range("K136").Value = 34280055
$ws.Range("L136").Value = 9498.999899999999
$ws.Range("M136").Value = -34277505
$ws.Range("N136").Value = -14598.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 25000
$ws.Range("J123").Value = 25000
$ws.Range("L123").Value = 25000
$ws.Range("N123").Value = -34800
# Row 135
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140
# Row 136
$ws.Range("H136").Value = 2015.2778
$ws.Range("I136").Value = 1572.3334
$ws.Range("J136").Value = 2901.1667
$ws.Range("K136").Value = 4717.0002
$ws.Range("L136").Value = 8703.500100000001
$ws.Range("M136").Value = -2167.0002
$ws.Range("N136").Value = -13803.5001
